$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.548.81'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '3.680.95'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '664.97'
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.77'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.497'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.145'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.10'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.69'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').Value = '3.678.25'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '69.532.30'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.118'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '16.04'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '467.23'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.74'
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.643'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '79.80'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '3.827.42'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000126'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.89'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.98'
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.66'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.68'
$ws.Range('E29').Value = '  -4.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.99'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.69'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.45'
$ws.Range('E33').Value = '  -3.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.163'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('D35').Value = '3.672.43'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.41'
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.04'
$ws.Range('E37').Value = '  -2.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '179.08'
$ws.Range('E39').Value = '  +4.62%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.21'
$ws.Range('E41').Value = '  -2.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0898'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.930'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '46.94'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.73'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.31'
$ws.Range('E47').Value = '  -3.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000269'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.06'
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.81'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('E51').Value = '  -1.88%  '
